$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 11 and 12 with corrected meter readings, and append
# rows 13-23 as new hourly readings for 2025-01-19 (today's DAP file).
#
# Columns: A=Date B=Hour C=S1 D=S2 E=S3 F=S4 G=S5 H=S6 I=S7 J=S8 K=Total
#          L=eS1 M=eS2 N=eS3 O=eS4 P=eS5 Q=eS6 R=eS7 S=eS8 T=eTotal

$rows = @(
    @{ Row=11; Date="2025-01-19"; Hour=10; C=20363; D=0; E=0; F=12258; G=0; H=0; I=0; J=0; K=32621; L=20365.0363; M=0; N=0; O=12259.2258; P=0; Q=0; R=0; S=0; T=32624.2621 },
    @{ Row=12; Date="2025-01-19"; Hour=11; C=20777; D=0; E=0; F=12410; G=0; H=0; I=0; J=0; K=33187; L=20779.0777; M=0; N=0; O=12411.241;  P=0; Q=0; R=0; S=0; T=33190.3187 },
    @{ Row=13; Date="2025-01-19"; Hour=12; C=20382; D=0; E=0; F=12615; G=0; H=0; I=0; J=0; K=32997; L=20384.0382; M=0; N=0; O=12616.2615; P=0; Q=0; R=0; S=0; T=33000.2997 },
    @{ Row=14; Date="2025-01-19"; Hour=13; C=21520; D=0; E=0; F=13127; G=0; H=0; I=0; J=0; K=34647; L=21522.152;  M=0; N=0; O=13128.3127; P=0; Q=0; R=0; S=0; T=34650.4647 },
    @{ Row=15; Date="2025-01-19"; Hour=14; C=22336; D=0; E=0; F=13457; G=0; H=0; I=0; J=0; K=35793; L=22338.2336; M=0; N=0; O=13458.3457; P=0; Q=0; R=0; S=0; T=35796.5793 },
    @{ Row=16; Date="2025-01-19"; Hour=15; C=22942; D=0; E=0; F=14118; G=0; H=0; I=0; J=0; K=37060; L=22944.2942; M=0; N=0; O=14119.4118; P=0; Q=0; R=0; S=0; T=37063.706  },
    @{ Row=17; Date="2025-01-19"; Hour=16; C=23137; D=0; E=0; F=14239; G=0; H=0; I=0; J=0; K=37376; L=23139.3137; M=0; N=0; O=14240.4239; P=0; Q=0; R=0; S=0; T=37379.7376 },
    @{ Row=18; Date="2025-01-19"; Hour=17; C=22951; D=0; E=0; F=14529; G=0; H=0; I=0; J=0; K=37480; L=22953.2951; M=0; N=0; O=14530.4529; P=0; Q=0; R=0; S=0; T=37483.748  },
    @{ Row=19; Date="2025-01-19"; Hour=18; C=24926; D=0; E=0; F=15330; G=0; H=0; I=0; J=0; K=40256; L=24928.4926; M=0; N=0; O=15331.533;  P=0; Q=0; R=0; S=0; T=40260.0256 },
    @{ Row=20; Date="2025-01-19"; Hour=19; C=25705; D=0; E=0; F=16529; G=0; H=0; I=0; J=0; K=42234; L=25707.5705; M=0; N=0; O=16530.6529; P=0; Q=0; R=0; S=0; T=42238.2234 },
    @{ Row=21; Date="2025-01-19"; Hour=20; C=24384; D=0; E=0; F=15922; G=0; H=0; I=0; J=0; K=40306; L=24386.4384; M=0; N=0; O=15923.5922; P=0; Q=0; R=0; S=0; T=40310.0306 },
    @{ Row=22; Date="2025-01-19"; Hour=21; C=24443; D=0; E=0; F=16382; G=0; H=0; I=0; J=0; K=40825; L=24445.4443; M=0; N=0; O=16383.6382; P=0; Q=0; R=0; S=0; T=40829.0825 },
    @{ Row=23; Date="2025-01-19"; Hour=22; C=24616; D=0; E=0; F=16531; G=0; H=0; I=0; J=0; K=41147; L=24618.4616; M=0; N=0; O=16532.6531; P=0; Q=0; R=0; S=0; T=41151.1147 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A holds a literal text date ("2025-01-19"), not a real date
    # serial. Force text formatting before the write so Excel doesn't
    # auto-convert the string into a date number, then drop back to the
    # Normal style so the new cell's <c> carries no style attribute.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $r.Hour
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
